# Weekly update: insert a new observation as row 219, pushing the
# existing rows 219:265 down to 220:266 (dimension grows to A1:R266).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 219 - shifts rows 219:265 down to 220:266
# and copies formatting from the row above (matches the s="2" date style
# already used on column D of the surrounding rows).
$ws.Rows.Item(219).Insert()

# Populate the new row 219 with this week's new record (same market /
# product as every other row in this single-category sheet).
$ws.Range("A219").Value = 9
$ws.Range("B219").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C219").Value = "Metropolitana"
$ws.Range("D219").Value = 44641
$ws.Range("E219").Value = 13
$ws.Range("F219").Value = 100112030
$ws.Range("G219").Value = "Poroto granado"
$ws.Range("H219").Value = "Sin especificar"
$ws.Range("I219").Value = "Primera"
$ws.Range("J219").Value = 25
$ws.Range("K219").Value = 24000
$ws.Range("L219").Value = 25000
$ws.Range("M219").Value = 24480
$ws.Range("N219").Value = "$/saco 25 kilos"
$ws.Range("O219").Value = "Provincia de Cardenal Caro"
$ws.Range("P219").Value = 979
$ws.Range("Q219").Value = 25
$ws.Range("R219").Value = "Hortaliza"
